# EPBDS: replaced bean factory configuration.
#
# The "classABeanFactory" column (column G) of the mapping table on
# Sheet1 - its header, its wrapped description ("Class A / bean factory /
# class") and every "org.dozer.factory.XMLBeanFactory" data value - is
# removed entirely. The following columns (oneWay / convertMethodAB)
# shift one column to the left (H->G, I->H).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the whole column; Excel shifts H/I left to G/H automatically.
$ws.Columns.Item(7).Delete() | Out-Null

# Row 17 no longer contains the wrapped "bean factory" header text, so it
# no longer needs the taller custom row height - let Excel recompute it.
$ws.Rows.Item(17).AutoFit() | Out-Null

# Park the selection where it ends up after the edit.
$ws.Range("G16").Select() | Out-Null
